$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bring D2,F2,G2,D3,F3,G3 in line with the rest of the quote-prefixed
# "text" cells in the block (copy the formatting only from B2, which is
# already set up that way).
$ws.Range("B2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("G3").PasteSpecial(-4122)

# Fix / normalise time + route values (remove stray colons from times,
# fix typo'd station code) and fill NULL data cells with a comma
# character so "\dep" checks downstream don't choke on blanks.
$ws.Range("D8").Value = "'1413"
$ws.Range("E8").Value = "'1403"
$ws.Range("F8").Value = "'1351,1357"
$ws.Range("G8").Value = "'1415"

$ws.Range("B4").Value = "'1627,1630,1703"

$ws.Range("G2").Value = "'KGX,LDS"

$ws.Range("C4").Value = "',"
$ws.Range("E4").Value = "',"
$ws.Range("F4").Value = "',"
$ws.Range("G4").Value = "',"
$ws.Range("C5").Value = "',"
$ws.Range("E5").Value = "',"
$ws.Range("F5").Value = "',"
$ws.Range("G5").Value = "',"
$ws.Range("C6").Value = "',"
$ws.Range("E6").Value = "',"
$ws.Range("F6").Value = "',"
$ws.Range("G6").Value = "',"
